$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.421.86"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.724.27"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'243.00"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "'0.4906"
$ws.Range("E7").Value = "  +1.97%  "
$ws.Range("D8").Value = "'0.2616"
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("D9").Value = "'0.06208"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "1.718.62"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "'0.07000"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").Value = "'4.564"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "'0.5992"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "'77.28"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "'0.9996"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "26.429.63"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "'0.9996"
$ws.Range("D19").Value = "'0.000007178"
$ws.Range("E19").Value = "  +3.07%  "
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("D21").Value = "1.939.82"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "'4.483"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("D24").Value = "'5.162"
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("D25").Value = "'138.17"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "'107.01"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "'1.713"
$ws.Range("E29").Value = "  -4.03%  "
$ws.Range("D30").Value = "'3.947"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").Value = "'0.07965"
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").Value = "'3.673"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").Value = "'0.04531"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.601"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'0.9948"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6261"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'0.9286"
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'1.959"
$ws.Range("E38").Value = "  -6.16%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.389"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "'0.9995"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.01486"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'99.59"
$ws.Range("E42").Value = "  -2.84%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.336"
$ws.Range("E43").Value = "  -3.59%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.3845"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'6.769"
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1168"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05369"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'7.737"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'30.12"
$ws.Range("E49").Value = "  -2.01%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.231"
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'50.86"
$ws.Range("E51").Value = "  -0.98%  "
